# Update test fixture Excel files for profile entities
# - Rename admin-level headers: County -> Province, Sub-County -> District, Ward -> Subdistrict
# - Make "School" (sheet 1) the active/selected tab instead of "Health Care Facilities"
# - Move the selection on the "School" sheet to G1

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)

# Row 1 headers on the "School" sheet: Name, Code, National, County, Sub-County, Ward, Village
$ws1.Range("D1").Value = "Province"
$ws1.Range("E1").Value = "District"
$ws1.Range("F1").Value = "Subdistrict"

# Switch the active sheet from "Health Care Facilities" back to "School",
# and move its selection to G1.
[void]$ws1.Activate()
[void]$ws1.Range("G1").Select()
